$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows right before the current row 355, shifting the
# existing rows 355-414 down to 361-420.
$ws.Range("A355:A360").EntireRow.Insert()

# Common (constant) column values shared by every data row in this sheet.
$marketId   = 11
$market     = "Vega Monumental Concepción"
$region     = "Bíobío"
$codreg     = 8
$tipo       = "Fruta"
$productId  = 100103
$producto   = "Frutos de hueso (carozo)"
$categoryId = 100103006
$categoria  = "Nectarín"
$origen     = "Región de O'Higgins"

# New weekly rows (fecha 44943) inserted at rows 355-360.
$rows = @(
    @{ Row=355; K="Sun Rise"; L="Especial"; M=50;  N=16000; O=16000; P=16000; Q="$/caja 15 kilos empedrada";   S=1067; T=15 },
    @{ Row=356; K="Sun Rise"; L="Primera";  M=50;  N=14000; O=14000; P=14000; Q="$/caja 15 kilos empedrada";   S=933;  T=15 },
    @{ Row=357; K="Sun Rise"; L="Segunda";  M=50;  N=12000; O=12000; P=12000; Q="$/caja 15 kilos empedrada";   S=800;  T=15 },
    @{ Row=358; K="Venus";    L="Especial"; M=50;  N=13000; O=13000; P=13000; Q="$/caja 15 kilos empedrada";   S=867;  T=15 },
    @{ Row=359; K="Venus";    L="Primera";  M=50;  N=11000; O=11000; P=11000; Q="$/caja 15 kilos empedrada";   S=733;  T=15 },
    @{ Row=360; K="Venus";    L="Segunda";  M=50;  N=10000; O=10000; P=10000; Q="$/caja 15 kilos empedrada";   S=667;  T=15 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value  = $marketId
    $ws.Cells.Item($rowNum, 2).Value  = $market
    $ws.Cells.Item($rowNum, 3).Value  = $region
    $ws.Cells.Item($rowNum, 4).Value  = 44943
    $ws.Cells.Item($rowNum, 5).Value  = $codreg
    $ws.Cells.Item($rowNum, 6).Value  = $tipo
    $ws.Cells.Item($rowNum, 7).Value  = $productId
    $ws.Cells.Item($rowNum, 8).Value  = $producto
    $ws.Cells.Item($rowNum, 9).Value  = $categoryId
    $ws.Cells.Item($rowNum, 10).Value = $categoria
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $origen
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
}
